$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (unchanged values, rewritten for clarity/completeness).
$ws.Range("B1").Value = "Win/Loss/Draw"
$ws.Range("C1").Value = "Role"

# Column B - Win/Loss/Draw result per row.
$ws.Range("B2").Value = "loss"
$ws.Range("B3").Value = "win"
$ws.Range("B4").Value = "win"
$ws.Range("B5").Value = "win"
$ws.Range("B6").Value = "win"
$ws.Range("B8").Value = "loss"

# Column C - new "Role" breakdown (tank specific wins/losses, plus the
# other roles) that didn't exist before this edit.
$ws.Range("C2").Value = "tank"
$ws.Range("C3").Value = "tank"
$ws.Range("C4").Value = "damage"
$ws.Range("C5").Value = "support"
$ws.Range("C6").Value = "tank"
$ws.Range("C8").Value = "damage"

# New rows 7 and 8, matching the border/number formatting already used by
# the existing A column cells (A2:A6).
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
